# Apply auto-update changes to the drone analysis sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (JOBY)
$ws.Range("D2").Value = 15.72
$ws.Range("E2").Value = 57.2
$ws.Range("F2").Value = 11.3
$ws.Range("K2").Value = 57.6
$ws.Range("N2").Value = 54.02451352198364

# Row 3 (ACHR)
$ws.Range("D3").Value = 8.73
$ws.Range("E3").Value = 59.3
$ws.Range("F3").Value = 16.49
$ws.Range("H3").Value = 40
$ws.Range("K3").Value = 54.2
$ws.Range("N3").Value = 54.02451352198364
